# Apply crypto price/volume updates per commit diff (Sun Dec 1 14:29:10 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "97.235.45"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "3.707.95"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'238.43"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("D7").Value = "'659.63"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("D8").Value = "'0.424"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").Value = "3.707.03"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "'0.0000321"
$ws.Range("E12").Value = "  +19.51%  "
$ws.Range("D13").Value = "'44.44"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "'6.84"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "4.399.25"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "96.938.57"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "'9.11"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "3.708.39"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'12.98"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'18.65"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  -4.46%  "
$ws.Range("D23").Value = "'520.97"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("D27").Value = "'102.26"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'0.196"
$ws.Range("E28").Value = "  +16.00%  "
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("D30").Value = "'12.82"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'656.09"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'32.28"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").Value = "'0.594"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'0.170"
$ws.Range("E41").Value = "  +4.91%  "
$ws.Range("D42").Value = "'6.84"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.486"
$ws.Range("E44").Value = "  +7.99%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'40.12"
$ws.Range("E45").Value = "  -8.55%  "
$ws.Range("D46").Value = "'0.970"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'0.0459"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "'2.33"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "'23.63"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'8.72"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("B51").Value = "MantraDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D51").Value = "'3.53"
$ws.Range("E51").Value = "  -3.83%  "
